$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Selection on the original sheet before switching away
[void]$ws1.Range("F10:F60").Select()

# Add the new "clean" sheet at the end of the workbook
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "clean"

# Headers (typed in this order: A, C, B, D, E)
$ws2.Range("A1").Value = "State"
$ws2.Range("C1").Value = "2021_pop"
$ws2.Range("B1").Value = "2020_pop"
$ws2.Range("D1").Value = "2022_pop"
$ws2.Range("E1").Value = "2023_pop"

# Copy per-state rows (State name + 2020-2023 population) from the main table
$ws1.Range("A10").Copy($ws2.Range("A2"))
$ws1.Range("C10").Copy($ws2.Range("B2"))
$ws1.Range("D10").Copy($ws2.Range("C2"))
$ws1.Range("E10").Copy($ws2.Range("D2"))
$ws1.Range("F10").Copy($ws2.Range("E2"))
$ws1.Range("A11").Copy($ws2.Range("A3"))
$ws1.Range("C11").Copy($ws2.Range("B3"))
$ws1.Range("D11").Copy($ws2.Range("C3"))
$ws1.Range("E11").Copy($ws2.Range("D3"))
$ws1.Range("F11").Copy($ws2.Range("E3"))
$ws1.Range("A12").Copy($ws2.Range("A4"))
$ws1.Range("C12").Copy($ws2.Range("B4"))
$ws1.Range("D12").Copy($ws2.Range("C4"))
$ws1.Range("E12").Copy($ws2.Range("D4"))
$ws1.Range("F12").Copy($ws2.Range("E4"))
$ws1.Range("A13").Copy($ws2.Range("A5"))
$ws1.Range("C13").Copy($ws2.Range("B5"))
$ws1.Range("D13").Copy($ws2.Range("C5"))
$ws1.Range("E13").Copy($ws2.Range("D5"))
$ws1.Range("F13").Copy($ws2.Range("E5"))
$ws1.Range("A14").Copy($ws2.Range("A6"))
$ws1.Range("C14").Copy($ws2.Range("B6"))
$ws1.Range("D14").Copy($ws2.Range("C6"))
$ws1.Range("E14").Copy($ws2.Range("D6"))
$ws1.Range("F14").Copy($ws2.Range("E6"))
$ws1.Range("A15").Copy($ws2.Range("A7"))
$ws1.Range("C15").Copy($ws2.Range("B7"))
$ws1.Range("D15").Copy($ws2.Range("C7"))
$ws1.Range("E15").Copy($ws2.Range("D7"))
$ws1.Range("F15").Copy($ws2.Range("E7"))
$ws1.Range("A16").Copy($ws2.Range("A8"))
$ws1.Range("C16").Copy($ws2.Range("B8"))
$ws1.Range("D16").Copy($ws2.Range("C8"))
$ws1.Range("E16").Copy($ws2.Range("D8"))
$ws1.Range("F16").Copy($ws2.Range("E8"))
$ws1.Range("A17").Copy($ws2.Range("A9"))
$ws1.Range("C17").Copy($ws2.Range("B9"))
$ws1.Range("D17").Copy($ws2.Range("C9"))
$ws1.Range("E17").Copy($ws2.Range("D9"))
$ws1.Range("F17").Copy($ws2.Range("E9"))
$ws1.Range("A18").Copy($ws2.Range("A10"))
$ws1.Range("C18").Copy($ws2.Range("B10"))
$ws1.Range("D18").Copy($ws2.Range("C10"))
$ws1.Range("E18").Copy($ws2.Range("D10"))
$ws1.Range("F18").Copy($ws2.Range("E10"))
$ws1.Range("A19").Copy($ws2.Range("A11"))
$ws1.Range("C19").Copy($ws2.Range("B11"))
$ws1.Range("D19").Copy($ws2.Range("C11"))
$ws1.Range("E19").Copy($ws2.Range("D11"))
$ws1.Range("F19").Copy($ws2.Range("E11"))
$ws1.Range("A20").Copy($ws2.Range("A12"))
$ws1.Range("C20").Copy($ws2.Range("B12"))
$ws1.Range("D20").Copy($ws2.Range("C12"))
$ws1.Range("E20").Copy($ws2.Range("D12"))
$ws1.Range("F20").Copy($ws2.Range("E12"))
$ws1.Range("A21").Copy($ws2.Range("A13"))
$ws1.Range("C21").Copy($ws2.Range("B13"))
$ws1.Range("D21").Copy($ws2.Range("C13"))
$ws1.Range("E21").Copy($ws2.Range("D13"))
$ws1.Range("F21").Copy($ws2.Range("E13"))
$ws1.Range("A22").Copy($ws2.Range("A14"))
$ws1.Range("C22").Copy($ws2.Range("B14"))
$ws1.Range("D22").Copy($ws2.Range("C14"))
$ws1.Range("E22").Copy($ws2.Range("D14"))
$ws1.Range("F22").Copy($ws2.Range("E14"))
$ws1.Range("A23").Copy($ws2.Range("A15"))
$ws1.Range("C23").Copy($ws2.Range("B15"))
$ws1.Range("D23").Copy($ws2.Range("C15"))
$ws1.Range("E23").Copy($ws2.Range("D15"))
$ws1.Range("F23").Copy($ws2.Range("E15"))
$ws1.Range("A24").Copy($ws2.Range("A16"))
$ws1.Range("C24").Copy($ws2.Range("B16"))
$ws1.Range("D24").Copy($ws2.Range("C16"))
$ws1.Range("E24").Copy($ws2.Range("D16"))
$ws1.Range("F24").Copy($ws2.Range("E16"))
$ws1.Range("A25").Copy($ws2.Range("A17"))
$ws1.Range("C25").Copy($ws2.Range("B17"))
$ws1.Range("D25").Copy($ws2.Range("C17"))
$ws1.Range("E25").Copy($ws2.Range("D17"))
$ws1.Range("F25").Copy($ws2.Range("E17"))
$ws1.Range("A26").Copy($ws2.Range("A18"))
$ws1.Range("C26").Copy($ws2.Range("B18"))
$ws1.Range("D26").Copy($ws2.Range("C18"))
$ws1.Range("E26").Copy($ws2.Range("D18"))
$ws1.Range("F26").Copy($ws2.Range("E18"))
$ws1.Range("A27").Copy($ws2.Range("A19"))
$ws1.Range("C27").Copy($ws2.Range("B19"))
$ws1.Range("D27").Copy($ws2.Range("C19"))
$ws1.Range("E27").Copy($ws2.Range("D19"))
$ws1.Range("F27").Copy($ws2.Range("E19"))
$ws1.Range("A28").Copy($ws2.Range("A20"))
$ws1.Range("C28").Copy($ws2.Range("B20"))
$ws1.Range("D28").Copy($ws2.Range("C20"))
$ws1.Range("E28").Copy($ws2.Range("D20"))
$ws1.Range("F28").Copy($ws2.Range("E20"))
$ws1.Range("A29").Copy($ws2.Range("A21"))
$ws1.Range("C29").Copy($ws2.Range("B21"))
$ws1.Range("D29").Copy($ws2.Range("C21"))
$ws1.Range("E29").Copy($ws2.Range("D21"))
$ws1.Range("F29").Copy($ws2.Range("E21"))
$ws1.Range("A30").Copy($ws2.Range("A22"))
$ws1.Range("C30").Copy($ws2.Range("B22"))
$ws1.Range("D30").Copy($ws2.Range("C22"))
$ws1.Range("E30").Copy($ws2.Range("D22"))
$ws1.Range("F30").Copy($ws2.Range("E22"))
$ws1.Range("A31").Copy($ws2.Range("A23"))
$ws1.Range("C31").Copy($ws2.Range("B23"))
$ws1.Range("D31").Copy($ws2.Range("C23"))
$ws1.Range("E31").Copy($ws2.Range("D23"))
$ws1.Range("F31").Copy($ws2.Range("E23"))
$ws1.Range("A32").Copy($ws2.Range("A24"))
$ws1.Range("C32").Copy($ws2.Range("B24"))
$ws1.Range("D32").Copy($ws2.Range("C24"))
$ws1.Range("E32").Copy($ws2.Range("D24"))
$ws1.Range("F32").Copy($ws2.Range("E24"))
$ws1.Range("A33").Copy($ws2.Range("A25"))
$ws1.Range("C33").Copy($ws2.Range("B25"))
$ws1.Range("D33").Copy($ws2.Range("C25"))
$ws1.Range("E33").Copy($ws2.Range("D25"))
$ws1.Range("F33").Copy($ws2.Range("E25"))
$ws1.Range("A34").Copy($ws2.Range("A26"))
$ws1.Range("C34").Copy($ws2.Range("B26"))
$ws1.Range("D34").Copy($ws2.Range("C26"))
$ws1.Range("E34").Copy($ws2.Range("D26"))
$ws1.Range("F34").Copy($ws2.Range("E26"))
$ws1.Range("A35").Copy($ws2.Range("A27"))
$ws1.Range("C35").Copy($ws2.Range("B27"))
$ws1.Range("D35").Copy($ws2.Range("C27"))
$ws1.Range("E35").Copy($ws2.Range("D27"))
$ws1.Range("F35").Copy($ws2.Range("E27"))
$ws1.Range("A36").Copy($ws2.Range("A28"))
$ws1.Range("C36").Copy($ws2.Range("B28"))
$ws1.Range("D36").Copy($ws2.Range("C28"))
$ws1.Range("E36").Copy($ws2.Range("D28"))
$ws1.Range("F36").Copy($ws2.Range("E28"))
$ws1.Range("A37").Copy($ws2.Range("A29"))
$ws1.Range("C37").Copy($ws2.Range("B29"))
$ws1.Range("D37").Copy($ws2.Range("C29"))
$ws1.Range("E37").Copy($ws2.Range("D29"))
$ws1.Range("F37").Copy($ws2.Range("E29"))
$ws1.Range("A38").Copy($ws2.Range("A30"))
$ws1.Range("C38").Copy($ws2.Range("B30"))
$ws1.Range("D38").Copy($ws2.Range("C30"))
$ws1.Range("E38").Copy($ws2.Range("D30"))
$ws1.Range("F38").Copy($ws2.Range("E30"))
$ws1.Range("A39").Copy($ws2.Range("A31"))
$ws1.Range("C39").Copy($ws2.Range("B31"))
$ws1.Range("D39").Copy($ws2.Range("C31"))
$ws1.Range("E39").Copy($ws2.Range("D31"))
$ws1.Range("F39").Copy($ws2.Range("E31"))
$ws1.Range("A40").Copy($ws2.Range("A32"))
$ws1.Range("C40").Copy($ws2.Range("B32"))
$ws1.Range("D40").Copy($ws2.Range("C32"))
$ws1.Range("E40").Copy($ws2.Range("D32"))
$ws1.Range("F40").Copy($ws2.Range("E32"))
$ws1.Range("A41").Copy($ws2.Range("A33"))
$ws1.Range("C41").Copy($ws2.Range("B33"))
$ws1.Range("D41").Copy($ws2.Range("C33"))
$ws1.Range("E41").Copy($ws2.Range("D33"))
$ws1.Range("F41").Copy($ws2.Range("E33"))
$ws1.Range("A42").Copy($ws2.Range("A34"))
$ws1.Range("C42").Copy($ws2.Range("B34"))
$ws1.Range("D42").Copy($ws2.Range("C34"))
$ws1.Range("E42").Copy($ws2.Range("D34"))
$ws1.Range("F42").Copy($ws2.Range("E34"))
$ws1.Range("A43").Copy($ws2.Range("A35"))
$ws1.Range("C43").Copy($ws2.Range("B35"))
$ws1.Range("D43").Copy($ws2.Range("C35"))
$ws1.Range("E43").Copy($ws2.Range("D35"))
$ws1.Range("F43").Copy($ws2.Range("E35"))
$ws1.Range("A44").Copy($ws2.Range("A36"))
$ws1.Range("C44").Copy($ws2.Range("B36"))
$ws1.Range("D44").Copy($ws2.Range("C36"))
$ws1.Range("E44").Copy($ws2.Range("D36"))
$ws1.Range("F44").Copy($ws2.Range("E36"))
$ws1.Range("A45").Copy($ws2.Range("A37"))
$ws1.Range("C45").Copy($ws2.Range("B37"))
$ws1.Range("D45").Copy($ws2.Range("C37"))
$ws1.Range("E45").Copy($ws2.Range("D37"))
$ws1.Range("F45").Copy($ws2.Range("E37"))
$ws1.Range("A46").Copy($ws2.Range("A38"))
$ws1.Range("C46").Copy($ws2.Range("B38"))
$ws1.Range("D46").Copy($ws2.Range("C38"))
$ws1.Range("E46").Copy($ws2.Range("D38"))
$ws1.Range("F46").Copy($ws2.Range("E38"))
$ws1.Range("A47").Copy($ws2.Range("A39"))
$ws1.Range("C47").Copy($ws2.Range("B39"))
$ws1.Range("D47").Copy($ws2.Range("C39"))
$ws1.Range("E47").Copy($ws2.Range("D39"))
$ws1.Range("F47").Copy($ws2.Range("E39"))
$ws1.Range("A48").Copy($ws2.Range("A40"))
$ws1.Range("C48").Copy($ws2.Range("B40"))
$ws1.Range("D48").Copy($ws2.Range("C40"))
$ws1.Range("E48").Copy($ws2.Range("D40"))
$ws1.Range("F48").Copy($ws2.Range("E40"))
$ws1.Range("A49").Copy($ws2.Range("A41"))
$ws1.Range("C49").Copy($ws2.Range("B41"))
$ws1.Range("D49").Copy($ws2.Range("C41"))
$ws1.Range("E49").Copy($ws2.Range("D41"))
$ws1.Range("F49").Copy($ws2.Range("E41"))
$ws1.Range("A50").Copy($ws2.Range("A42"))
$ws1.Range("C50").Copy($ws2.Range("B42"))
$ws1.Range("D50").Copy($ws2.Range("C42"))
$ws1.Range("E50").Copy($ws2.Range("D42"))
$ws1.Range("F50").Copy($ws2.Range("E42"))
$ws1.Range("A51").Copy($ws2.Range("A43"))
$ws1.Range("C51").Copy($ws2.Range("B43"))
$ws1.Range("D51").Copy($ws2.Range("C43"))
$ws1.Range("E51").Copy($ws2.Range("D43"))
$ws1.Range("F51").Copy($ws2.Range("E43"))
$ws1.Range("A52").Copy($ws2.Range("A44"))
$ws1.Range("C52").Copy($ws2.Range("B44"))
$ws1.Range("D52").Copy($ws2.Range("C44"))
$ws1.Range("E52").Copy($ws2.Range("D44"))
$ws1.Range("F52").Copy($ws2.Range("E44"))
$ws1.Range("A53").Copy($ws2.Range("A45"))
$ws1.Range("C53").Copy($ws2.Range("B45"))
$ws1.Range("D53").Copy($ws2.Range("C45"))
$ws1.Range("E53").Copy($ws2.Range("D45"))
$ws1.Range("F53").Copy($ws2.Range("E45"))
$ws1.Range("A54").Copy($ws2.Range("A46"))
$ws1.Range("C54").Copy($ws2.Range("B46"))
$ws1.Range("D54").Copy($ws2.Range("C46"))
$ws1.Range("E54").Copy($ws2.Range("D46"))
$ws1.Range("F54").Copy($ws2.Range("E46"))
$ws1.Range("A55").Copy($ws2.Range("A47"))
$ws1.Range("C55").Copy($ws2.Range("B47"))
$ws1.Range("D55").Copy($ws2.Range("C47"))
$ws1.Range("E55").Copy($ws2.Range("D47"))
$ws1.Range("F55").Copy($ws2.Range("E47"))
$ws1.Range("A56").Copy($ws2.Range("A48"))
$ws1.Range("C56").Copy($ws2.Range("B48"))
$ws1.Range("D56").Copy($ws2.Range("C48"))
$ws1.Range("E56").Copy($ws2.Range("D48"))
$ws1.Range("F56").Copy($ws2.Range("E48"))
$ws1.Range("A57").Copy($ws2.Range("A49"))
$ws1.Range("C57").Copy($ws2.Range("B49"))
$ws1.Range("D57").Copy($ws2.Range("C49"))
$ws1.Range("E57").Copy($ws2.Range("D49"))
$ws1.Range("F57").Copy($ws2.Range("E49"))
$ws1.Range("A58").Copy($ws2.Range("A50"))
$ws1.Range("C58").Copy($ws2.Range("B50"))
$ws1.Range("D58").Copy($ws2.Range("C50"))
$ws1.Range("E58").Copy($ws2.Range("D50"))
$ws1.Range("F58").Copy($ws2.Range("E50"))
$ws1.Range("A59").Copy($ws2.Range("A51"))
$ws1.Range("C59").Copy($ws2.Range("B51"))
$ws1.Range("D59").Copy($ws2.Range("C51"))
$ws1.Range("E59").Copy($ws2.Range("D51"))
$ws1.Range("F59").Copy($ws2.Range("E51"))
$ws1.Range("A60").Copy($ws2.Range("A52"))
$ws1.Range("C60").Copy($ws2.Range("B52"))
$ws1.Range("D60").Copy($ws2.Range("C52"))
$ws1.Range("E60").Copy($ws2.Range("D52"))
$ws1.Range("F60").Copy($ws2.Range("E52"))

# Final selection/active-cell on the new sheet
[void]$ws2.Range("G50").Select()
